# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund holding detail, same shape as the
# existing quarterly sheets) right before the "总计" summary sheet, and adds
# a leading "2022-Q1" row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before "总计".
#    (Worksheet handles returned by Worksheets.Item(...) track *position*,
#    not identity, in this engine - re-fetch "总计" by name after the sheet
#    is inserted, since the just-created sheet now occupies the old slot.)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. Header row + the A-column "index" style - copied from the "2021-Q4"
#    sheet so the new sheet picks up the same cell style (s="2": bold,
#    bordered, centered) the other quarterly detail sheets use.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$refSheet.Range("B1:H1").Copy()
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

$refSheet2 = $wb.Worksheets.Item("2021-Q4")
$refSheet2.Range("A2").Copy()
$q1Sheet2 = $wb.Worksheets.Item("2022-Q1")
$q1Sheet2.Range("A2").PasteSpecial(-4122)

$ws = $wb.Worksheets.Item("2022-Q1")
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Detail rows. Columns B (fund code) and D:G (scale/position/value,
#    stored as text in every other quarterly sheet) are typed as text so
#    leading zeros and exact decimal text are preserved; the NumberFormat
#    is reset to "Normal" afterwards so no stray style index lingers on the
#    cell (matches the plain, un-styled cells used elsewhere in the file).
# ---------------------------------------------------------------------------
$rows = @(
    @(0,  "513330", "华夏恒生互联网科技业ETF（QDII）", "233.65", "96.98", "3.62", "8.4581", 9),
    @(1,  "513180", "华夏恒生科技交易型开放式指数证券投资基金（QDII）", "89.43", "94.63", "3.91", "3.4967", 10),
    @(2,  "513010", "易方达恒生科技交易型开放式指数证券投资基金（QDII）", "26.24", "94.46", "3.87", "1.0155", 10),
    @(3,  "159740", "大成恒生科技交易型开放式指数证券投资基金（QDII）", "7.06", "98.90", "4.06", "0.2866", 10),
    @(4,  "159742", "博时恒生科技交易型开放式指数证券投资基金(QDII)", "3.73", "96.57", "3.96", "0.1477", 10),
    @(5,  "159741", "嘉实恒生科技交易型开放式指数证券投资基金（QDII）", "3.29", "99.73", "4.09", "0.1346", 10),
    @(6,  "501021", "华宝兴业标普香港上市中国中小盘指数(QDII-LOF)A", "4.96", "94.77", "2.01", "0.0997", 6),
    @(7,  "513890", "上投摩根恒生科技ETF（QDII）", "2.10", "93.73", "3.85", "0.0808", 10),
    @(8,  "013127", "汇添富恒生科技指数（QDII）A", "1.51", "91.71", "3.76", "0.0568", 10),
    @(9,  "513160", "银华恒生港股通中国科技ETF", "0.62", "92.07", "7.23", "0.0448", 7),
    @(10, "013128", "汇添富恒生科技指数（QDII）C", "0.63", "91.71", "3.76", "0.0237", 10),
    @(11, "006127", "华宝兴业标普香港上市中国中小盘指数(QDII-LOF)C", "0.23", "94.77", "2.01", "0.0046", 6),
    @(12, "005269", "华泰柏瑞港股通量化灵活配置混合", "0.33", "37.77", "0.87", "0.0029", 8),
    @(13, "004321", "前海开源沪港深强国产业灵活配置混合", "0.02", "64.32", "4.63", "0.0009", 3)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]

    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("B$r").Style = "Normal"

    $ws.Range("C$r").Value = $row[2]

    $ws.Range("D$r`:G$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("D$r`:G$r").Style = "Normal"

    $ws.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. Update the "总计" summary sheet: insert a new leading data row for
#    2022-Q1 (shifting the existing quarters down) and renumber the index
#    column.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 14
$totalSheet.Range("D2").Value = 13.85

# Renumber the index column (A) for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
